$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 619.5
$ws.Range("I96").Value = 450.84616
$ws.Range("K96").Value = 1352.53848
$ws.Range("M96").Value = 20.46152000000006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 900
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 598
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 900
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -250
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 40524.668
$ws.Range("J128").Value = 40524.668
$ws.Range("L128").Value = 40524.668
$ws.Range("N128").Value = -50484.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1581.125
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1581.125
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 4743.375
$ws.Range("N137").Value = -9843.375
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5298.75
$ws.Range("I138").Value = 5600
$ws.Range("J138").Value = 5263.7207
$ws.Range("K138").Value = 16800
$ws.Range("L138").Value = 15791.1621
$ws.Range("M138").Value = -11660
$ws.Range("N138").Value = -26071.1621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19657.596
$ws.Range("I32").Value = 3558.4126
$ws.Range("K32").Value = 3558.4126
$ws.Range("M32").Value = -3271.4126

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1814.2609
$ws.Range("I61").Value = 886.4666999999999
$ws.Range("J61").Value = 2263.1936
$ws.Range("K61").Value = 886.4666999999999
$ws.Range("L61").Value = 2263.1936
$ws.Range("M61").Value = -674.4666999999999
$ws.Range("N61").Value = -2687.1936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 18561.2
$ws.Range("J80").Value = 18561.2
$ws.Range("L80").Value = 18561.2
$ws.Range("N80").Value = -20557.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 18561.2
$ws.Range("J83").Value = 18561.2
$ws.Range("L83").Value = 55683.60000000001
$ws.Range("N83").Value = -65667.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3908
$ws.Range("I132").Value = 5285.3335
$ws.Range("J132").Value = 2100.25
$ws.Range("K132").Value = 15856.0005
$ws.Range("L132").Value = 6300.75
$ws.Range("M132").Value = -13326.0005
$ws.Range("N132").Value = -11360.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1814.2609
$ws.Range("I136").Value = 886.4666999999999
$ws.Range("J136").Value = 2263.1936
$ws.Range("K136").Value = 2659.4001
$ws.Range("L136").Value = 6789.5808
$ws.Range("M136").Value = -109.4000999999998
$ws.Range("N136").Value = -11889.5808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1496.3334
$ws.Range("I99").Value = 986.6667
$ws.Range("J99").Value = 2006
$ws.Range("K99").Value = 986.6667
$ws.Range("L99").Value = 2006
$ws.Range("M99").Value = 511.3333
$ws.Range("N99").Value = -5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 15876119
$ws.Range("I107").Value = 19610418
$ws.Range("J107").Value = 5344.5
$ws.Range("K107").Value = 19610418
$ws.Range("L107").Value = 5344.5
$ws.Range("M107").Value = -19608498
$ws.Range("N107").Value = -9184.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15844.681
$ws.Range("I31").Value = 38280.89
$ws.Range("J31").Value = 2382.9556
$ws.Range("K31").Value = 38280.89
$ws.Range("L31").Value = 2382.9556
$ws.Range("M31").Value = -37985.89
$ws.Range("N31").Value = -2972.9556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 15844.681
$ws.Range("I34").Value = 38280.89
$ws.Range("J34").Value = 2382.9556
$ws.Range("K34").Value = 38280.89
$ws.Range("L34").Value = 2382.9556
$ws.Range("M34").Value = -38078.89
$ws.Range("N34").Value = -2786.9556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1235.6428
$ws.Range("J94").Value = 1316.909
$ws.Range("L94").Value = 1316.909
$ws.Range("N94").Value = -2218.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 725.2
$ws.Range("I107").Value = 747.61536
$ws.Range("J107").Value = 700.9167
$ws.Range("K107").Value = 747.61536
$ws.Range("L107").Value = 700.9167
$ws.Range("M107").Value = 1172.38464
$ws.Range("N107").Value = -4540.9167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 68000
$ws.Range("J140").Value = 68000
$ws.Range("L140").Value = 68000
$ws.Range("N140").Value = -78360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1291.5
$ws.Range("I34").Value = 116.666664
$ws.Range("J34").Value = 2466.3333
$ws.Range("K34").Value = 349.999992
$ws.Range("L34").Value = 7398.999899999999
$ws.Range("M34").Value = -265.999992
$ws.Range("N34").Value = -7566.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2126.3484
$ws.Range("I68").Value = 1404.3928
$ws.Range("J68").Value = 2658.3157
$ws.Range("K68").Value = 4213.178400000001
$ws.Range("L68").Value = 7974.9471
$ws.Range("M68").Value = -3402.178400000001
$ws.Range("N68").Value = -9596.947100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2126.3484
$ws.Range("I71").Value = 1404.3928
$ws.Range("J71").Value = 2658.3157
$ws.Range("K71").Value = 12639.5352
$ws.Range("L71").Value = 23924.8413
$ws.Range("M71").Value = -8583.5352
$ws.Range("N71").Value = -32036.8413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6512.52
$ws.Range("I75").Value = 471
$ws.Range("J75").Value = 7336.364
$ws.Range("K75").Value = 1413
$ws.Range("L75").Value = 22009.092
$ws.Range("M75").Value = -415
$ws.Range("N75").Value = -24005.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 6512.52
$ws.Range("I78").Value = 471
$ws.Range("J78").Value = 7336.364
$ws.Range("K78").Value = 4239
$ws.Range("L78").Value = 66027.276
$ws.Range("M78").Value = 753
$ws.Range("N78").Value = -76011.276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1755933.5
$ws.Range("I131").Value = 587.8
$ws.Range("J131").Value = 2021895
$ws.Range("K131").Value = 1763.4
$ws.Range("L131").Value = 6065685
$ws.Range("M131").Value = 3276.6
$ws.Range("N131").Value = -6075765

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2451.6667
$ws.Range("I132").Value = 1770.4
$ws.Range("K132").Value = 5311.200000000001
$ws.Range("M132").Value = -2781.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7038.4614
$ws.Range("I132").Value = 10500.5
$ws.Range("J132").Value = 4071
$ws.Range("K132").Value = 31501.5
$ws.Range("L132").Value = 12213
$ws.Range("M132").Value = -28971.5
$ws.Range("N132").Value = -17273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 20715
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 20715
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 20715
$ws.Range("N139").Value = -30995
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2363.2
$ws.Range("I122").Value = 2238.6667
$ws.Range("J122").Value = 2550
$ws.Range("K122").Value = 6716.000100000001
$ws.Range("L122").Value = 7650
$ws.Range("M122").Value = -4266.000100000001
$ws.Range("N122").Value = -12550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3170
$ws.Range("I132").Value = 3934.35
$ws.Range("J132").Value = 2078.0715
$ws.Range("K132").Value = 11803.05
$ws.Range("L132").Value = 6234.2145
$ws.Range("M132").Value = -9273.049999999999
$ws.Range("N132").Value = -11294.2145
